$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "product backlog"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("product backlog")

# Row 20 becomes a "highlighted" (yellow, wrap-text) row like the other
# finished backlog entries, and gets a Status of "(3) Completed".
$ws1.Range("A20:I20").Interior.Color = 65535
$ws1.Range("A20:I20").WrapText = $true
$ws1.Range("I20").Value = "(3) Completed"

# Row 21 gains a note in column E describing the new backlog item.
$ws1.Range("E21").Value = "enable oath 2 for authorization"

# ---------------------------------------------------------------------
# Sheet "sprint backlog"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("sprint backlog")

$ws2.Range("A2").Value = 11
$ws2.Range("E2").Value = "to see pdf reports for employees per department etc."

$ws2.Range("A3").Value = 21
$ws2.Range("E3").Value = "enable oath 2 for authorization"
$ws2.Rows(3).AutoFit()

$ws2.Range("I1:I1048576").Validation.Delete()
$ws2.Range("I1:I1048576").Validation.Add(3, 1, 1, "userstorystatus")

$ws2.Range("H1:H1048576").Validation.Delete()
$ws2.Range("H1:H1048576").Validation.Add(3, 1, 1, "priority")

# ---------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------
$ws1.Range("A21:XFD21").Select()

$ws2.Range("E6").Select()
